$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "settings" sheet: remove the form_id column (column B).
#    Column C (version) -> B, D (style) -> C, E (namespaces) -> D.
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

# Comments are anchored to their cell and do not travel with a column
# delete, so copy the text left by hand before removing the column.
$versionComment    = $settings.Range("C1").Comment.Text()
$styleComment      = $settings.Range("D1").Comment.Text()
$namespacesComment = $settings.Range("E1").Comment.Text()

[void]$settings.Range("B1").Comment.Text($versionComment)
[void]$settings.Range("C1").Comment.Text($styleComment)
[void]$settings.Range("D1").Comment.Text($namespacesComment)
$settings.Range("E1").Comment.Delete()

# Now actually remove the form_id column; this shifts the remaining
# columns (values + shared strings) left automatically.
$settings.Columns("B").Delete()

# Keep the recorded selection consistent with the shifted columns
# (old selection was on the now-removed E11, which lines up with C11).
# Selecting a range activates its sheet, so re-activate "survey"
# afterwards to keep the originally active tab unchanged.
[void]$settings.Range("C11").Select()
[void]$wb.Worksheets.Item("survey").Activate()

# ---------------------------------------------------------------------------
# 2) "survey" sheet: tidy up the fragmented conditional-formatting ranges
#    that had been split around row 27 (no real data lives there) back
#    into simple, uniform ranges.
# ---------------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")
$fcs = $survey.Cells.FormatConditions

# First rule group (A28:H10000 A27:B27 D27:H27 A2:H26) -> A2:H10000
$fcs.Item(1).ModifyAppliesToRange($survey.Range("A2:H10000"))

# C-column rule group (C28:C10000 C2:C26) -> C2:C10000
$fcs.Item(7).ModifyAppliesToRange($survey.Range("C2:C10000"))

# Drop the extra duplicate rule groups that only applied to C27
for ($i = 14; $i -ge 9; $i--) {
    $fcs.Item($i).Delete()
}
